# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to reflect the latest generated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 86
    $ws.Range("F4").Value = 2251
}
